$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.966.49"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.635.75"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.55"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.72%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.54"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.84%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.864.18"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.648.67"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.19"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.990.24"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₃0745"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.75"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.34%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "191.00"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.25"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.64"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.13"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.90%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.44"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.71%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.96%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.28"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.66%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.59%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.39%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.84%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.32%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.50"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.46%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.41"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.137.78"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.866"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.42%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.524"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.19%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.60"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.778"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.56%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.773.32"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.63%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.28"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0529"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.43%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.57"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.65%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.06%  "
